$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 1 de Abril de 2020 a las 05:20'

$ws.Range('B4').Value = 188578
$ws.Range('C4').Value = 48
$ws.Range('E4').Value = 177437
$ws.Range('G4').Value = 1
$ws.Range('H4').Value = 3890

$ws.Range('B22').Value = 4804
$ws.Range('C22').Value = 41
$ws.Range('E22').Value = 4439

$ws.Range('A79').Value = 'Kazajistan'
$ws.Range('B79').Value = 348
$ws.Range('C79').Value = 5
$ws.Range('D79').Value = 24
$ws.Range('E79').Value = 322
$ws.Range('F79').Value = 6

$ws.Range('A80').Value = 'Costa Rica'
$ws.Range('B80').Value = 347
$ws.Range('D80').Value = 4
$ws.Range('E80').Value = 341
$ws.Range('F80').Value = 8

$ws.Range('A100').Value = 'Honduras'
$ws.Range('B100').Value = 172
$ws.Range('C100').Value = 31
$ws.Range('D100').Value = 3
$ws.Range('E100').Value = 159
$ws.Range('F100').Value = 4
$ws.Range('G100').Value = 3
$ws.Range('H100').Value = 10

$ws.Range('A101').Value = 'Malta'
$ws.Range('D101').Value = 2
$ws.Range('E101').Value = 167
$ws.Range('F101').Value = 2

$ws.Range('A102').Value = 'Islas Feroe'
$ws.Range('B102').Value = 169
$ws.Range('D102').Value = 74
$ws.Range('E102').Value = 95
$ws.Range('F102').Value = 3
$ws.Range('H102').Value = 0

$ws.Range('A103').Value = 'Ghana'
$ws.Range('B103').Value = 161
$ws.Range('D103').Value = 31
$ws.Range('E103').Value = 125
$ws.Range('F103').Value = 1
$ws.Range('H103').Value = 5

$ws.Range('A104').Value = 'Bielorrusia'
$ws.Range('B104').Value = 152
$ws.Range('D104').Value = 47
$ws.Range('E104').Value = 104
$ws.Range('F104').Value = 2
$ws.Range('H104').Value = 1

$ws.Range('A105').Value = 'Mauricio'
$ws.Range('D105').Value = 0
$ws.Range('E105').Value = 138
$ws.Range('F105').Value = 1
$ws.Range('H105').Value = 5

$ws.Range('A106').Value = 'Sri Lanka'
$ws.Range('B106').Value = 143
$ws.Range('D106').Value = 17
$ws.Range('E106').Value = 124
$ws.Range('F106').Value = 5
$ws.Range('H106').Value = 2

$ws.Range('E120').Value = 82
$ws.Range('G120').Value = 1
$ws.Range('H120').Value = 4

$ws.Range('A142').Value = 'El Salvador'
$ws.Range('F142').Value = 5

$ws.Range('A143').Value = 'Guam'
$ws.Range('F143').Value = 0

$ws.Range('A156').Value = 'Bahamas'
$ws.Range('C156').Value = 1

$ws.Range('A159').Value = 'Guinea Ecuatorial'
$ws.Range('D159').Value = 1
$ws.Range('E159').Value = 14
$ws.Range('H159').Value = 0

$ws.Range('A160').Value = 'San Martin (Parte Francesa)'
$ws.Range('B160').Value = 15
$ws.Range('D160').Value = 2
$ws.Range('E160').Value = 12

$ws.Range('A161').Value = 'Islas Caimanes'
$ws.Range('D161').Value = 0
$ws.Range('H161').Value = 1

$ws.Range('A168').Value = 'Seychelles'
$ws.Range('C168').Value = 0

$ws.Range('A169').Value = 'Surinam'
$ws.Range('C169').Value = 1

$ws.Range('A173').Value = 'Laos'

$ws.Range('A175').Value = 'Granada'

$ws.Range('A177').Value = 'Guinea-Bisau'

$ws.Range('A178').Value = 'San Cristobal y Nieves'

$ws.Range('A179').Value = 'Mozambique'

$ws.Range('A190').Value = 'Islas Turcas y Caicos'

$ws.Range('A191').Value = 'Fiyi'

$ws.Range('A192').Value = 'Montserrat'

$ws.Range('A193').Value = 'Nicaragua'
$ws.Range('D193').Value = 0
$ws.Range('H193').Value = 1

$ws.Range('A195').Value = 'Somalia'
$ws.Range('D195').Value = 1
$ws.Range('H195').Value = 0

$ws.Range('A203').Value = 'Anguila'

$ws.Range('A204').Value = 'Bonaire, San Eustaquio y Saba'
$ws.Range('C204').Value = 2

$ws.Range('A205').Value = 'Burundi'
$ws.Range('B205').Value = 2
$ws.Range('E205').Value = 2

$ws.Range('A206').Value = 'Timor Oriental'

$ws.Range('A208').Value = 'Sierra Leona'
$ws.Range('D208').Value = 0
$ws.Range('E208').Value = 1

$ws.Range('A209').Value = 'San Vicente y las Granadinas'
$ws.Range('B209').Value = 1
$ws.Range('C209').Value = 0
$ws.Range('D209').Value = 1
$ws.Range('E209').Value = 0
$ws.Range('F209').Value = 0
$ws.Range('G209').Value = 0
$ws.Range('H209').Value = 0
